$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-20"

# Update the 2022 column header text
$ws.Range("I1").Value = "2022 (through 10-20)"

# Update October 2022 value
$ws.Range("I11").Value = 67

# Update Total 2022 value
$ws.Range("I14").Value = 1344
